# TIVA based Daughter Manual Board First Draft
# Adds two new small "Pin / Plug and Play / uC based" tables in columns L:N
# next to the existing "White Line Sensors" (rows 3-5) and "Sharp Sensors"
# (rows 8-10) tables, describing UART and I2C connector pinouts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- UART pinout table (rows 3-5), values entered in the order the author
# typed them so new shared-string entries land in the matching slots.
$ws.Range("M3").Value = "Plug and Play"
$ws.Range("N3").Value = "uC based"
$ws.Range("L3").Value = "Pin"

$ws.Range("L4").Value = "RX"
$ws.Range("L5").Value = "TX"

$ws.Range("M4").Value = "PB0(UART1)"
$ws.Range("M5").Value = "PB1(UART1)"

$ws.Range("N4").Value = "PC6(UART3)"
$ws.Range("N5").Value = "PC7(UART3)"

# --- I2C pinout table (rows 8-10) ---
$ws.Range("L8").Value = "Pin"
$ws.Range("M8").Value = "Plug and Play"
$ws.Range("N8").Value = "uC based"

$ws.Range("L9").Value = "SCL"
$ws.Range("L10").Value = "SDA"

$ws.Range("N10").Value = "PB3(I2C0)"
$ws.Range("N9").Value = "PB2(I2C0)"

$ws.Range("M9").Value = "PA6(I2C1)"
$ws.Range("M10").Value = "PA7(I2C1)"

# --- Formatting: bold headers, column widths ---
$ws.Range("L3:N3").Font.Bold = $true
$ws.Range("L8:N8").Font.Bold = $true

$ws.Columns.Item(13).ColumnWidth = 13.166666666666666
$ws.Columns.Item(14).ColumnWidth = 10.5

# Update selection to match the final cursor position left by the author
$ws.Range("M15").Select()
